$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N (14th column) on the
# "Repayment schedule" sheet - this shifts the old N/O/P ("Late",
# "heading"/"Date", "Outstanding") columns one position to the right
# (-> O/P/Q), matching the "Variable Instalments" layout change.
$ws.Columns.Item(14).Insert()

# Match the new column's width to the neighbouring columns (~10.71 chars).
$ws.Columns.Item(14).ColumnWidth = 9.8

# Make "Repayment schedule" the active sheet / tab, with cell K14 selected
# (this also clears the previous tab-selection on "Transactions").
$ws.Activate() | Out-Null
$ws.Range("K14").Select() | Out-Null
